# Auto-generated Excel COM-interop script to apply cryptos.xlsx diff
# Commit: Updated symbol list on Thu Jan 19 05:12:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All cells in columns D, E and G on this sheet are stored as plain text
# (inline/shared strings), even though many look numeric or percentage-like.
# Force text format before assigning so Excel does not silently convert
# these numeric-looking strings into actual numbers.

# --- Column G: uniform change, "4" -> "5" for every data row (2-51) ---
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "5"

# --- Column D: updated Price values ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "290.97"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.80"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.938"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07187"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.778"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.661"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8948"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1653"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07708"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08057"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03026"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1001"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001490"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005740"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.469"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04511"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004009"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01597"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04379"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007376"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007688"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002005"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009236"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005994"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002998"

# --- Column E: updated Volume(1h) percentage values ---
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.00%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.92%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.14%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-8.18%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-12.46%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.22%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.91%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.92%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.26%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.55%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.61%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.31%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.43%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.88%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.58%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.01%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.89%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.38%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-6.42%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.37%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.34%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.95%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-9.98%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.06%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-8.17%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-8.55%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.83%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.24%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-15.04%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.48%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.13%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.13%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "172.72%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.33%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.13%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.13%"
